# Add a new "metadata" worksheet after the existing "data" sheet and
# populate it with the panel query metadata (data_name, data_id,
# data_version, data_version_created, panel_query_time, panel_get_request).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Insert the new sheet right after "data" so the tab order is data, metadata.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$ws.Name = "metadata"

# Header row (row 1) - columns B..G.
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Match the bold/bordered/centered header style already used on the
# "data" sheet's header row, by copying its formatting over.
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Data row (row 2).
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Testicular cancer pertinent cancer susceptibility"
$ws.Range("C2").Value = 82

# Force "1.1" to be stored as text (not the number 1.1), then drop the
# temporary text number-format so the cell keeps the default style.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.1"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").Value = "2020-05-07T14:26:17.781188Z"
$ws.Range("F2").Value = "2021-10-05 14:22:54.454605"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/82/?format=json"

# A2 also carries the header-style formatting (bold/bordered/centered).
$dataSheet.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

Write-Output "metadata sheet added"
